$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("I17").Value = "sd"
$ws.Range("J17").Value = "Statement-non-opinion"
$ws.Range("I19").Value = "%"
$ws.Range("J19").Value = "Uninterpretable"
$ws.Range("I23").Value = "sd"
$ws.Range("J23").Value = "Statement-non-opinion"
$ws.Range("I37").Value = "sd"
$ws.Range("J37").Value = "Statement-non-opinion"
$ws.Range("I39").Value = "sv"
$ws.Range("J39").Value = "Statement-opinion"
$ws.Range("I52").Value = "aa"
$ws.Range("J52").Value = "Agree/Accept"
$ws.Range("I63").Value = "aa"
$ws.Range("J63").Value = "Agree/Accept"
$ws.Range("I88").Value = "ba"
$ws.Range("J88").Value = "Appreciation"
$ws.Range("I90").Value = "sd"
$ws.Range("J90").Value = "Statement-non-opinion"
$ws.Range("I100").Value = "sv"
$ws.Range("J100").Value = "Statement-opinion"
$ws.Range("I102").Value = "sd"
$ws.Range("J102").Value = "Statement-non-opinion"
$ws.Range("I112").Value = "sd"
$ws.Range("J112").Value = "Statement-non-opinion"
$ws.Range("I118").Value = "sv"
$ws.Range("J118").Value = "Statement-opinion"
$ws.Range("I130").Value = "sd"
$ws.Range("J130").Value = "Statement-non-opinion"
$ws.Range("I137").Value = "sd"
$ws.Range("J137").Value = "Statement-non-opinion"
$ws.Range("I175").Value = "sd"
$ws.Range("J175").Value = "Statement-non-opinion"
$ws.Range("I176").Value = "sd"
$ws.Range("J176").Value = "Statement-non-opinion"
$ws.Range("I192").Value = "sv"
$ws.Range("J192").Value = "Statement-opinion"
$ws.Range("I202").Value = "sv"
$ws.Range("J202").Value = "Statement-opinion"
$ws.Range("I206").Value = "ba"
$ws.Range("J206").Value = "Appreciation"
$ws.Range("I216").Value = "sv"
$ws.Range("J216").Value = "Statement-opinion"
$ws.Range("I218").Value = "aa"
$ws.Range("J218").Value = "Agree/Accept"
$ws.Range("I244").Value = "sv"
$ws.Range("J244").Value = "Statement-opinion"
$ws.Range("I257").Value = "aa"
$ws.Range("J257").Value = "Agree/Accept"
$ws.Range("I261").Value = "sv"
$ws.Range("J261").Value = "Statement-opinion"
$ws.Range("I265").Value = "sd"
$ws.Range("J265").Value = "Statement-non-opinion"
$ws.Range("I266").Value = "sv"
$ws.Range("J266").Value = "Statement-opinion"
$ws.Range("I270").Value = "aa"
$ws.Range("J270").Value = "Agree/Accept"
$ws.Range("I272").Value = "sd"
$ws.Range("J272").Value = "Statement-non-opinion"
$ws.Range("I284").Value = "sd"
$ws.Range("J284").Value = "Statement-non-opinion"
$ws.Range("I295").Value = "sd"
$ws.Range("J295").Value = "Statement-non-opinion"
$ws.Range("I310").Value = "sd"
$ws.Range("J310").Value = "Statement-non-opinion"
$ws.Range("I311").Value = "sd"
$ws.Range("J311").Value = "Statement-non-opinion"
$ws.Range("I314").Value = "sv"
$ws.Range("J314").Value = "Statement-opinion"
$ws.Range("I320").Value = "sd"
$ws.Range("J320").Value = "Statement-non-opinion"
$ws.Range("I328").Value = "sv"
$ws.Range("J328").Value = "Statement-opinion"
$ws.Range("I329").Value = "sv"
$ws.Range("J329").Value = "Statement-opinion"
$ws.Range("I341").Value = "sd"
$ws.Range("J341").Value = "Statement-non-opinion"
$ws.Range("I345").Value = "sd"
$ws.Range("J345").Value = "Statement-non-opinion"
